$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 20
$ws.Range("B2").Value = 84
$ws.Range("B3").Value = 136
$ws.Range("B4").Value = 194
$ws.Range("B5").Value = 225
$ws.Range("B6").Value = 247
$ws.Range("B7").Value = 269
$ws.Range("B8").Value = 289
$ws.Range("B9").Value = 314
